# RPA datasets push 2024-05-18
# Insert two new IPO rows at the top of the data (row 2), pushing the
# existing data rows down by two. Matches the canonical diff: rows for
# "KB제28호스팩" and "아이씨티케이" are prepended ahead of "코칩".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for two new data rows right after the header row. Excel's
# Insert() copies the formatting of the row above (the bold header), so
# strip it back to the plain/default style used by the rest of the data
# rows before writing any values into them.
$ws.Range("A2:A3").EntireRow.Insert()
$ws.Range("A2:Y3").Style = "Normal"

# Date-like and percent-like text ("2024-04-29", "6.54%", ...) must stay
# literal strings (as in every other data row), not get auto-converted by
# Excel into a date serial / percentage number. Force Text format first,
# then restore the plain "Normal" style afterwards so the saved cell has
# no explicit format, matching the rest of the sheet.
$dateRange = $ws.Range("A2:C3")
$dateRange.NumberFormat = "@"

$pctRange = $ws.Range("O2:O3")
$pctRange.NumberFormat = "@"

# New row 2: KB제28호스팩
$ws.Range("A2").Value = "2024-04-29"
$ws.Range("B2").Value = "2024-04-30"
$ws.Range("C2").Value = "2024-05-17"
$ws.Range("D2").Value = "KB"
$ws.Range("E2").Value = "KB제28호스팩"
$ws.Range("F2").Value = 5000000
$ws.Range("G2").Value = 5000000
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 2000
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 5505000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 2000
$ws.Range("N2").Value = "1118.39:1"
$ws.Range("O2").Value = "-"
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = "기업인수합병"

# New row 3: 아이씨티케이
$ws.Range("A3").Value = "2024-04-24"
$ws.Range("B3").Value = "2024-05-30"
$ws.Range("C3").Value = "2024-05-17"
$ws.Range("D3").Value = "NH"
$ws.Range("E3").Value = "아이씨티케이"
$ws.Range("F3").Value = 1970000
$ws.Range("G3").Value = 1970000
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 13000
$ws.Range("J3").Value = 16000
$ws.Range("K3").Value = 13124496
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 20000
$ws.Range("N3").Value = "783.2:1"
$ws.Range("O3").Value = "6.54%"
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 0
$ws.Range("Y3").Value = "PUF반도체,보안솔루션(보안반도체,정보통신모듈기기,정보통신용반도체) 제조,개발"

# Restore the default (unstyled) appearance now that the literal text is
# locked in, so the new rows match the formatting of the rest of the table.
$ws.Range("A2:Y3").Style = "Normal"
